# Weekly update: insert two new price records at the top of the data
# block (rows 293-294), pushing all existing records down by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 293 (each Insert() shifts rows 293.. down by one).
$ws.Rows.Item(293).Insert()
$ws.Rows.Item(293).Insert()

# New row 293: Primera, 180 cajas, $/caja 12 unidades, Ecuador.
$ws.Cells.Item(293,1).Value  = 10
$ws.Cells.Item(293,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(293,3).Value  = "La Araucanía"
$ws.Cells.Item(293,4).Value  = 44588
$ws.Cells.Item(293,5).Value  = 9
$ws.Cells.Item(293,6).Value  = "Fruta"
$ws.Cells.Item(293,7).Value  = 100108
$ws.Cells.Item(293,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(293,9).Value  = 100108005
$ws.Cells.Item(293,10).Value = "Piña"
$ws.Cells.Item(293,11).Value = "Caramelo"
$ws.Cells.Item(293,12).Value = "Primera"
$ws.Cells.Item(293,13).Value = 180
$ws.Cells.Item(293,14).Value = 18000
$ws.Cells.Item(293,15).Value = 19000
$ws.Cells.Item(293,16).Value = 18444
$ws.Cells.Item(293,17).Value = "$/caja 12 unidades"
$ws.Cells.Item(293,18).Value = "Ecuador"
$ws.Cells.Item(293,19).Value = 1537
$ws.Cells.Item(293,20).Value = 12

# New row 294: Segunda, 100 cajas, $/caja 14 unidades, Ecuador.
$ws.Cells.Item(294,1).Value  = 10
$ws.Cells.Item(294,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(294,3).Value  = "La Araucanía"
$ws.Cells.Item(294,4).Value  = 44588
$ws.Cells.Item(294,5).Value  = 9
$ws.Cells.Item(294,6).Value  = "Fruta"
$ws.Cells.Item(294,7).Value  = 100108
$ws.Cells.Item(294,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(294,9).Value  = 100108005
$ws.Cells.Item(294,10).Value = "Piña"
$ws.Cells.Item(294,11).Value = "Caramelo"
$ws.Cells.Item(294,12).Value = "Segunda"
$ws.Cells.Item(294,13).Value = 100
$ws.Cells.Item(294,14).Value = 18000
$ws.Cells.Item(294,15).Value = 18000
$ws.Cells.Item(294,16).Value = 18000
$ws.Cells.Item(294,17).Value = "$/caja 14 unidades"
$ws.Cells.Item(294,18).Value = "Ecuador"
$ws.Cells.Item(294,19).Value = 1286
$ws.Cells.Item(294,20).Value = 14
